$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for existing rows 334-364 (columns F and G)
$updates = @{
    334 = @{ F = 204779; G = 3482 }
    335 = @{ F = 130403; G = 2975 }
    336 = @{ F = 101007; G = 3296 }
    337 = @{ F = 103273; G = 2882 }
    338 = @{ F = 224324; G = 3126 }
    339 = @{ F = 653549; G = 5583 }
    340 = @{ F = 380355; G = 3256 }
    341 = @{ F = 294314; G = 3649 }
    342 = @{ F = 179512; G = 3087 }
    343 = @{ F = 133568; G = 2958 }
    344 = @{ F = 136445; G = 2526 }
    345 = @{ F = 289977; G = 3301 }
    346 = @{ F = 658647; G = 4692 }
    347 = @{ F = 338827; G = 2876 }
    348 = @{ F = 230879; G = 3222 }
    349 = @{ F = 159515 }
    350 = @{ F = 128655; G = 2778 }
    351 = @{ F = 149543; G = 2803 }
    352 = @{ F = 304940; G = 3529 }
    353 = @{ F = 712622; G = 5223 }
    354 = @{ F = 303431; G = 2768 }
    355 = @{ F = 221026; G = 3410 }
    356 = @{ F = 159723; G = 2875 }
    357 = @{ F = 137945; G = 3008 }
    358 = @{ F = 160771 }
    359 = @{ F = 318660; G = 3329 }
    360 = @{ F = 728145; G = 4950 }
    361 = @{ F = 325693; G = 2556 }
    362 = @{ F = 219272; G = 3022 }
    363 = @{ F = 179839; G = 2663 }
    364 = @{ F = 156538; G = 2301 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('F')) {
        $ws.Range("F$row").Value = $vals['F']
    }
    if ($vals.ContainsKey('G')) {
        $ws.Range("G$row").Value = $vals['G']
    }
}

# Add new row 365
$ws.Range("A365").Value = 44259
$ws.Range("B365").Value = 319582
$ws.Range("C365").Value = 11749
$ws.Range("D365").Value = 2423
$ws.Range("E365").Value = 7665
$ws.Range("F365").Value = 151964
$ws.Range("G365").Value = 2027
